{"js": "// Add a centered, underlined title paragraph (\"Technical Writing Cheat\n// sheet\") above the existing intro paragraph, and re-flow the intro\n// paragraph's runs (merging \"H\" + \"ere is ...\" + \" I have gather ... writer\"\n// + \":\" into three runs, with a grammar-check proofErr pair wrapped around\n// the repeated word \"Technical\").\n\nconst body = context.document.body;\nconst firstPara = body.paragraphs.getFirst();\nawait context.sync();\n\n// --- 1. Insert the new title paragraph before the current first paragraph.\n// insertParagraph() inherits the anchor paragraph's pPr/rPr (Title style,\n// 360 line spacing, Times New Roman 12pt run), so we only need to layer the\n// title-specific formatting (centered, 13pt/26 half-points, underline) on\n// top of it.\nconst titlePara = firstPara.insertParagraph(\n  \"Technical Writing Cheat sheet\",\n  Word.InsertLocation.before\n);\ntitlePara.alignment = Word.Alignment.centered;\ntitlePara.font.size = 13;\ntitlePara.font.sizeBidirectional = 13;\ntitlePara.font.underline = Word.UnderlineType.single;\nawait context.sync();\n\n// Also stamp the paragraph-mark (end-of-paragraph) run properties so the\n// pPr/rPr matches the run rPr exactly, same as Word does when you type a\n// whole paragraph in that format.\nconst titleEnd = titlePara.getRange(Word.RangeLocation.end);\ntitleEnd.font.size = 13;\ntitleEnd.font.sizeBidirectional = 13;\ntitleEnd.font.underline = Word.UnderlineType.single;\nawait context.sync();\n\n// --- 2. Re-flow the runs of the original intro paragraph (now the second\n// paragraph in the body). Its paragraph properties are unchanged; only the\n// run/text structure changes, gaining a proofErr gramStart/gramEnd pair\n// around the (now duplicated) word \"Technical\".\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst introPara = paragraphs.items[1];\nconst introRange = introPara.getRange();\n\nconst introOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr>' +\n  '<w:pStyle w:val=\"Title\"/>' +\n  '<w:spacing w:line=\"360\" w:lineRule=\"auto\"/>' +\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:sz w:val=\"24\"/>' +\n  '<w:szCs w:val=\"24\"/>' +\n  '</w:rPr>' +\n  '</w:pPr>' +\n  '<w:r>' +\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:sz w:val=\"24\"/>' +\n  '<w:szCs w:val=\"24\"/>' +\n  '</w:rPr>' +\n  '<w:t xml:space=\"preserve\">Here is a technical writing cheat sheet that covers some important concepts and best practices I have gather along the why as a </w:t>' +\n  '</w:r>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r>' +\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:sz w:val=\"24\"/>' +\n  '<w:szCs w:val=\"24\"/>' +\n  '</w:rPr>' +\n  '<w:t>Technical</w:t>' +\n  '</w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r>' +\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>' +\n  '<w:sz w:val=\"24\"/>' +\n  '<w:szCs w:val=\"24\"/>' +\n  '</w:rPr>' +\n  '<w:t xml:space=\"preserve\"> writer:</w:t>' +\n  '</w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nintroRange.insertOoxml(introOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Add a centered, underlined title paragraph (\"Technical Writing Cheat\n# sheet\") above the existing intro paragraph, and re-flow the intro\n# paragraph's runs (merging \"H\" + \"ere is ...\" + \" I have gather ... writer\"\n# + \":\" into three runs, with a grammar-check proofErr pair wrapped around\n# the repeated word \"Technical\").\n\n$d = $word.ActiveDocument\n\n# --- 1. Insert the new title paragraph before the current first paragraph.\n# InsertParagraphBefore() inherits the anchor paragraph's pPr/rPr (Title\n# style, 360 line spacing, Times New Roman 12pt run), so we only need to\n# layer the title-specific formatting (centered, 13pt/26 half-points,\n# underline) on top of it.\n$firstPara = $d.Paragraphs(1)\n$firstPara.Range.InsertParagraphBefore()\n\n$titlePara = $d.Paragraphs(1)\n$titleRng = $titlePara.Range\n$titleRng.Text = \"Technical Writing Cheat sheet\"\n$titlePara.Alignment = 1\n$titleRng.Font.Size = 13\n$titleRng.Font.SizeBi = 13\n$titleRng.Font.Underline = 1\n\n# --- 2. Re-flow the runs of the original intro paragraph (now the second\n# paragraph in the document). Its paragraph properties are unchanged; only\n# the run/text structure changes, gaining a proofErr gramStart/gramEnd pair\n# around the (now duplicated) word \"Technical\".\n$introPara = $d.Paragraphs(2)\n$introRng = $introPara.Range\n\n$introXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"Title\"/><w:spacing w:line=\"360\" w:lineRule=\"auto\"/><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">Here is a technical writing cheat sheet that covers some important concepts and best practices I have gather along the why as a </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Technical</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> writer:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$introRng.InsertXML($introXml)\n"}
